# Automatische test-sync: 2025-06-22 19:02:50
# Append the new "B2B samenwerkingsvoorstel" log entry to the Logs sheet,
# extend the conditional formatting ranges to include the new row, and
# refresh the category-count roll-up on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 32 with the new mail-log entry.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(32, 1).Value = "B2B samenwerkingsvoorstel"
$logs.Cells.Item(32, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(32, 3).Value = "Ik wil graag een samenwerking bespreken voor onze zakelijke klanten."
$logs.Cells.Item(32, 4).Value = "Samenwerking / Partnerverzoek"
$logs.Cells.Item(32, 5).Value = "Geachte heer/mevrouw,`nHartelijk dank voor uw interesse in een mogelijke samenwerking met ons bedrijf voor onze zakelijke klanten. Wij waarderen uw voorstel en willen graag meer informatie ontvangen om te bekijken hoe we kunnen samenwerken.`nKunt u meer details verstrekken over uw bedrijf en de diensten die u aanbiedt? Op basis hiervan kunnen we beoordelen of er mogelijkheden zijn voor een vruchtbare samenwerking.`nIk kijk uit naar uw antwoord.`nMet vriendelijke groet,`n[Naam]  `nE-mailassistent  `n[Bedrijfsnaam]"
$logs.Cells.Item(32, 6).Value = "2025-06-22 19:02:34"
$logs.Cells.Item(32, 7).Value = "Ja"

# The multi-line reply in column E makes the COM layer stamp an explicit
# (wrap-driven) row height on row 32; AutoFit() resolves it back to the
# sheet's computed/standard height so no stray customHeight sticks around,
# matching the other multi-line rows already in the sheet.
$logs.Rows.Item(32).AutoFit()

# Extend the two conditional-formatting blocks (Categorie / Beantwoord
# columns) so they keep covering the full data range through row 32.
$logs.Range("D2:D31").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D32"))
$logs.Range("G2:G31").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G32"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: Samenwerking / Partnerverzoek now has 4 mails
#    (tied with IT / Technisch probleem and Retour / Terugbetaling), so
#    it moves up above Productinformatie / Afmelding in the ranking.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Samenwerking / Partnerverzoek"
$dash.Cells.Item(3, 2).Value = 4

$dash.Cells.Item(5, 1).Value = "Productinformatie"
$dash.Cells.Item(5, 2).Value = 4

$dash.Cells.Item(6, 1).Value = "Afmelding / Nieuwsbrief"
$dash.Cells.Item(6, 2).Value = 3
